$d = $word.ActiveDocument

# Locate the title-page line that reports the sheet count, e.g. "Листов 12".
$para = $d.Content
$null = $para.Find.Execute("Листов", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
$para.Expand(4)

# Within that paragraph, find the "12" run precisely.
$target = $d.Range($para.Start, $para.End)
$null = $target.Find.Execute("12", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)

# $target now spans exactly "12". Change only its second character ("2" -> "1"),
# toggling a character-formatting property around the edit so the engine keeps
# the untouched first character in the original run and places the edited
# character into a new run of its own (matching how Word splits runs when a
# single character inside a larger run is edited in place).
$secondChar = $d.Range($target.End - 1, $target.End)
$secondChar.Font.Bold = $true
$secondChar.Text = "1"
$secondCharAfter = $d.Range($target.End - 1, $target.End)
$secondCharAfter.Font.Bold = $false
